$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 22.0484878138354
$ws.Range("C2").Value = 7.086576677024713
$ws.Range("D2").Value = 9.772756460511282
$ws.Range("E2").Value = 10.02013192933944
$ws.Range("F2").Value = 51.48609614823092
$ws.Range("L2").Value = 10.38116335073162
$ws.Range("B3").Value = 21.76138917492624
$ws.Range("C3").Value = 6.627830347426846
$ws.Range("D3").Value = 9.654196310775246
$ws.Range("E3").Value = 10.01792632901276
$ws.Range("F3").Value = 50.33935074908573
$ws.Range("L3").Value = 10.37744982034082
$ws.Range("B4").Value = 21.59331492376979
$ws.Range("C4").Value = 6.329235690193566
$ws.Range("D4").Value = 9.581002176143727
$ws.Range("E4").Value = 10.01709169825472
$ws.Range("F4").Value = 49.63020320537586
$ws.Range("L4").Value = 10.37764010772753
$ws.Range("B5").Value = 21.52698135237814
$ws.Range("C5").Value = 6.20329161062159
$ws.Range("D5").Value = 9.551089646392899
$ws.Range("E5").Value = 10.01688176157242
$ws.Range("F5").Value = 49.34028504999038
$ws.Range("L5").Value = 10.37833839773096
$ws.Range("B6").Value = 21.51609971108232
$ws.Range("C6").Value = 6.18212090492607
$ws.Range("D6").Value = 9.546117965559333
$ws.Range("E6").Value = 10.01685475019199
$ws.Range("F6").Value = 49.29209781294129
$ws.Range("L6").Value = 10.37849181933987
$ws.Range("B7").Value = 21.59241146710756
$ws.Range("C7").Value = 6.327554441328208
$ws.Range("D7").Value = 9.580599091536566
$ws.Range("E7").Value = 10.0170883405126
$ws.Range("F7").Value = 49.62629660586742
$ws.Range("L7").Value = 10.37764701258604
$ws.Range("B8").Value = 21.94785164635859
$ws.Range("C8").Value = 6.931907290069574
$ws.Range("D8").Value = 9.731968318559609
$ws.Range("E8").Value = 10.01926330444802
$ws.Range("F8").Value = 51.09194122699942
$ws.Range("L8").Value = 10.37937001928289
$ws.Range("B9").Value = 22.70556283401936
$ws.Range("C9").Value = 7.983254905177745
$ws.Range("D9").Value = 10.02500069203515
$ws.Range("E9").Value = 10.02767349806251
$ws.Range("F9").Value = 53.91201034155129
$ws.Range("L9").Value = 10.4023569889279
$ws.Range("B10").Value = 23.29308280119702
$ws.Range("C10").Value = 8.674878739815513
$ws.Range("D10").Value = 10.23711248379156
$ws.Range("E10").Value = 10.0364109055182
$ws.Range("F10").Value = 55.93310628958004
$ws.Range("L10").Value = 10.43118803550282
$ws.Range("B11").Value = 23.56570214304283
$ws.Range("C11").Value = 8.97214989870403
$ws.Range("D11").Value = 10.33273607416806
$ws.Range("E11").Value = 10.04094794405009
$ws.Range("F11").Value = 56.83805544552168
$ws.Range("L11").Value = 10.44688716085988
$ws.Range("B12").Value = 23.66959241502713
$ws.Range("C12").Value = 9.082238192240068
$ws.Range("D12").Value = 10.3688064050319
$ws.Range("E12").Value = 10.04274735343877
$ws.Range("F12").Value = 57.17839089387282
$ws.Range("L12").Value = 10.45320200225992
$ws.Range("B13").Value = 23.64719045526257
$ws.Range("C13").Value = 9.058638872222541
$ws.Range("D13").Value = 10.3610444959619
$ws.Range("E13").Value = 10.04235619390334
$ws.Range("F13").Value = 57.10520215084141
$ws.Range("L13").Value = 10.45182556552832
$ws.Range("B14").Value = 23.57423677825729
$ws.Range("C14").Value = 8.981256590749652
$ws.Range("D14").Value = 10.33570650078839
$ws.Range("E14").Value = 10.04109434990626
$ws.Range("F14").Value = 56.86610361032442
$ws.Range("L14").Value = 10.44739928346159
$ws.Range("B15").Value = 23.52963249567581
$ws.Range("C15").Value = 8.933534903139813
$ws.Range("D15").Value = 10.32016752619054
$ws.Range("E15").Value = 10.04033203942851
$ws.Range("F15").Value = 56.71933521663861
$ws.Range("L15").Value = 10.44473617834693
$ws.Range("B16").Value = 23.27536439424644
$ws.Range("C16").Value = 8.655103335785752
$ws.Range("D16").Value = 10.23084447217567
$ws.Range("E16").Value = 10.03612575538822
$ws.Range("F16").Value = 55.87365129645875
$ws.Range("L16").Value = 10.43021392287683
$ws.Range("B17").Value = 23.12066850303977
$ws.Range("C17").Value = 8.479859689487673
$ws.Range("D17").Value = 10.17581442349057
$ws.Range("E17").Value = 10.03368970264779
$ws.Range("F17").Value = 55.35095063802551
$ws.Range("L17").Value = 10.42196571995097
$ws.Range("B18").Value = 23.03220224888766
$ws.Range("C18").Value = 8.377430098040168
$ws.Range("D18").Value = 10.14408189171612
$ws.Range("E18").Value = 10.0323414532526
$ws.Range("F18").Value = 55.0489612424597
$ws.Range("L18").Value = 10.41746488673992
$ws.Range("B19").Value = 23.00234032116628
$ws.Range("C19").Value = 8.342468143715905
$ws.Range("D19").Value = 10.13332441950534
$ws.Range("E19").Value = 10.03189403105685
$ws.Range("F19").Value = 54.94649021206052
$ws.Range("L19").Value = 10.41598281440405
$ws.Range("B20").Value = 23.13708411956565
$ws.Range("C20").Value = 8.49868370053963
$ws.Range("D20").Value = 10.18168093219141
$ws.Range("E20").Value = 10.03394354587338
$ws.Range("F20").Value = 55.40673431633124
$ws.Range("L20").Value = 10.42281858354283
$ws.Range("B21").Value = 23.59564816571287
$ws.Range("C21").Value = 9.004052890659722
$ws.Range("D21").Value = 10.34315281294269
$ws.Range("E21").Value = 10.04146277310713
$ws.Range("F21").Value = 56.93639846475143
$ws.Range("L21").Value = 10.44868936495741
$ws.Range("B22").Value = 23.89911344564233
$ws.Range("C22").Value = 9.319881258566564
$ws.Range("D22").Value = 10.44785956249896
$ws.Range("E22").Value = 10.04685129296852
$ws.Range("F22").Value = 57.92231747493458
$ws.Range("L22").Value = 10.46775259402379
$ws.Range("B23").Value = 23.73684111789482
$ws.Range("C23").Value = 9.15263625717583
$ws.Range("D23").Value = 10.3920560960543
$ws.Range("E23").Value = 10.04393179150389
$ws.Range("F23").Value = 57.39746157877413
$ws.Range("L23").Value = 10.45738162325292
$ws.Range("B24").Value = 23.12966114520587
$ws.Range("C24").Value = 8.490178597992506
$ws.Range("D24").Value = 10.17902897840229
$ws.Range("E24").Value = 10.03382862057745
$ws.Range("F24").Value = 55.38151911022066
$ws.Range("L24").Value = 10.42243225260891
$ws.Range("B25").Value = 22.4947113268589
$ws.Range("C25").Value = 7.713085235959547
$ws.Range("D25").Value = 9.946223733472165
$ws.Range("E25").Value = 10.02495125711182
$ws.Range("F25").Value = 53.15679566544176
$ws.Range("L25").Value = 10.39403885515932